$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shorten column B (KOR headword) cells down to the bare term --
# the longer fragments were leftover context from an earlier extraction pass.
$ws.Range("B4").Value = '척수내압'
$ws.Range("B7").Value = '경막천자'
$ws.Range("B9").Value = '척수관류압'
$ws.Range("B10").Value = '비결핵항산균'
$ws.Range("B12").Value = '말 검사'
$ws.Range("B18").Value = '직장암'
$ws.Range("B28").Value = '분절하기관지'
$ws.Range("B33").Value = '의 핵산길잡이'
$ws.Range("B45").Value = '안정액'
$ws.Range("B48").Value = '척추마취'
$ws.Range("B49").Value = '신경아교증'
$ws.Range("B53").Value = '혁신공간'
$ws.Range("B54").Value = '합의공간'
$ws.Range("B61").Value = '규제약물법'
$ws.Range("B65").Value = '금융'
$ws.Range("B68").Value = '빅 블러'
$ws.Range("B71").Value = '레그테크'
$ws.Range("B73").Value = '큐'
$ws.Range("B74").Value = '고관여'
$ws.Range("B75").Value = '저관여'
$ws.Range("B80").Value = '관광'
$ws.Range("B89").Value = '크리스터  페터손'
$ws.Range("B93").Value = '소프트스킬'
$ws.Range("B96").Value = '태평양 경제 사회'
$ws.Range("B99").Value = '트위터'
$ws.Range("B104").Value = '입력 스트림'
$ws.Range("B109").Value = '관행 구'
$ws.Range("B117").Value = '유럽 개인정보보호법'
$ws.Range("B118").Value = '버트랜드'
$ws.Range("B122").Value = '적층 제조'
$ws.Range("B123").Value = '기술'
$ws.Range("B126").Value = '다우'
$ws.Range("B130").Value = '혁신공간'
$ws.Range("B131").Value = '합의공간'
$ws.Range("B147").Value = '남용'
$ws.Range("B148").Value = '전신마취'

# B87 held a stray leftover fragment with no matching English term; clear it.
$ws.Range("B87").ClearContents()

# New row of data appended at the bottom of the table.
$ws.Range("A150").Value = '이 연구는 크레인의 조작장치 상의 휴먼에러 도 포함하고 있으나 설계, 제작, 보전, 운전, 폐기 등크레인과 관련된 전 수명주기(life cycle)적인 내용을 담고 있어 조작장치의 양립성 등 인간공학적 문제에 집중된 결과를 내지 못했다.'
$ws.Range("B150").Value = '전, 운전, 폐기'
$ws.Range("C150").Value = 'life cycle'
